# Scheduled-runner refresh of market-price-derived columns (H-N) across
# several sheets: currentAveragePrice / currentAveragePriceNQ / HQ,
# LevePriceNQ / HQ, and the derived LeveProfitNQ / HQ figures.
# Values below are the refreshed snapshot; a couple of rows gain/lose
# the M/N profit cells entirely where the underlying price data now
# does/doesn't clear the leve cost threshold.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3449.75
$ws.Range("I64").Value = 3174.5
$ws.Range("J64").Value = 3725
$ws.Range("K64").Value = 3174.5
$ws.Range("L64").Value = 3725
$ws.Range("M64").Value = -2926.5
$ws.Range("N64").Value = -4221
$ws.Range("H67").Value = 3449.75
$ws.Range("I67").Value = 3174.5
$ws.Range("J67").Value = 3725
$ws.Range("K67").Value = 3174.5
$ws.Range("L67").Value = 3725
$ws.Range("M67").Value = -2316.5
$ws.Range("N67").Value = -5441

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2369.9092
$ws.Range("I2").Value = 2385.6
$ws.Range("J2").Value = 2213
$ws.Range("K2").Value = 2385.6
$ws.Range("L2").Value = 2213
$ws.Range("M2").Value = -2272.6
$ws.Range("N2").Value = -2439
$ws.Range("H23").Value = 15181.875
$ws.Range("J23").Value = 15181.875
$ws.Range("L23").Value = 15181.875
$ws.Range("N23").Value = -15699.875
$ws.Range("H32").Value = 17549486
$ws.Range("I32").Value = 29414036
$ws.Range("J32").Value = 10586.305
$ws.Range("K32").Value = 29414036
$ws.Range("L32").Value = 10586.305
$ws.Range("M32").Value = -29413749
$ws.Range("N32").Value = -11160.305
$ws.Range("H45").Value = 2595.9473
$ws.Range("I45").Value = 2380.9375
$ws.Range("J45").Value = 3742.6667
$ws.Range("K45").Value = 2380.9375
$ws.Range("L45").Value = 3742.6667
$ws.Range("M45").Value = -2003.9375
$ws.Range("N45").Value = -4496.6667
$ws.Range("H61").Value = 2892.9473
$ws.Range("I61").Value = 2023.8
$ws.Range("J61").Value = 3858.6667
$ws.Range("K61").Value = 2023.8
$ws.Range("L61").Value = 3858.6667
$ws.Range("M61").Value = -1811.8
$ws.Range("N61").Value = -4282.6667
$ws.Range("H74").Value = 3091.2917
$ws.Range("I74").Value = 3841.75
$ws.Range("J74").Value = 1590.375
$ws.Range("K74").Value = 3841.75
$ws.Range("L74").Value = 1590.375
$ws.Range("M74").Value = -2967.75
$ws.Range("N74").Value = -3338.375
$ws.Range("H77").Value = 3091.2917
$ws.Range("I77").Value = 3841.75
$ws.Range("J77").Value = 1590.375
$ws.Range("K77").Value = 19208.75
$ws.Range("L77").Value = 7951.875
$ws.Range("M77").Value = -14840.75
$ws.Range("N77").Value = -16687.875
$ws.Range("H116").Value = 2369.9092
$ws.Range("I116").Value = 2385.6
$ws.Range("J116").Value = 2213
$ws.Range("K116").Value = 2385.6
$ws.Range("L116").Value = 2213
$ws.Range("M116").Value = -91.59999999999991
$ws.Range("N116").Value = -6801
$ws.Range("H132").Value = 3756
$ws.Range("I132").Value = 3340.7693
$ws.Range("J132").Value = 4655.6665
$ws.Range("K132").Value = 10022.3079
$ws.Range("L132").Value = 13966.9995
$ws.Range("M132").Value = -7492.3079
$ws.Range("N132").Value = -19026.9995
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H136").Value = 2892.9473
$ws.Range("I136").Value = 2023.8
$ws.Range("J136").Value = 3858.6667
$ws.Range("K136").Value = 6071.4
$ws.Range("L136").Value = 11576.0001
$ws.Range("M136").Value = -3521.4
$ws.Range("N136").Value = -16676.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2369.9092
$ws.Range("I3").Value = 2385.6
$ws.Range("J3").Value = 2213
$ws.Range("K3").Value = 2385.6
$ws.Range("L3").Value = 2213
$ws.Range("M3").Value = -2271.6
$ws.Range("N3").Value = -2441
$ws.Range("H64").Value = 528.8
$ws.Range("I64").Value = 967.8
$ws.Range("K64").Value = 967.8
$ws.Range("M64").Value = -742.8
$ws.Range("H67").Value = 528.8
$ws.Range("I67").Value = 967.8
$ws.Range("K67").Value = 967.8
$ws.Range("M67").Value = -187.8
$ws.Range("H86").Value = 2538.7273
$ws.Range("I86").Value = 2082.4
$ws.Range("J86").Value = 2919
$ws.Range("K86").Value = 2082.4
$ws.Range("L86").Value = 2919
$ws.Range("M86").Value = -959.4000000000001
$ws.Range("N86").Value = -5165
$ws.Range("H89").Value = 2538.7273
$ws.Range("I89").Value = 2082.4
$ws.Range("J89").Value = 2919
$ws.Range("K89").Value = 10412
$ws.Range("L89").Value = 14595
$ws.Range("M89").Value = -4796
$ws.Range("N89").Value = -25827
$ws.Range("H94").Value = 752.1539
$ws.Range("I94").Value = 659.75
$ws.Range("J94").Value = 900
$ws.Range("K94").Value = 659.75
$ws.Range("L94").Value = 900
$ws.Range("M94").Value = -208.75
$ws.Range("N94").Value = -1802
$ws.Range("H122").Value = 49491.668
$ws.Range("J122").Value = 49491.668
$ws.Range("L122").Value = 49491.668
$ws.Range("N122").Value = -59291.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H63").Value = 23000
$ws.Range("I63").Value = 21000
$ws.Range("J63").Value = 25000
$ws.Range("K63").Value = 21000
$ws.Range("L63").Value = 25000
$ws.Range("M63").Value = -20314
$ws.Range("N63").Value = -26372
$ws.Range("H66").Value = 23000
$ws.Range("I66").Value = 21000
$ws.Range("J66").Value = 25000
$ws.Range("K66").Value = 63000
$ws.Range("L66").Value = 75000
$ws.Range("M66").Value = -59568
$ws.Range("N66").Value = -81864
$ws.Range("H94").Value = 3223.0908
$ws.Range("I94").Value = 1314.25
$ws.Range("J94").Value = 4313.857
$ws.Range("K94").Value = 1314.25
$ws.Range("L94").Value = 4313.857
$ws.Range("M94").Value = -863.25
$ws.Range("N94").Value = -5215.857
$ws.Range("H132").Value = 50002204
$ws.Range("I132").Value = 71430010
$ws.Range("J132").Value = 3987.6667
$ws.Range("K132").Value = 214290030
$ws.Range("L132").Value = 11963.0001
$ws.Range("M132").Value = -214287500
$ws.Range("N132").Value = -17023.0001
$ws.Range("H134").Value = 3109.4883
$ws.Range("I134").Value = 1740.7587
$ws.Range("J134").Value = 5944.7144
$ws.Range("K134").Value = 5222.2761
$ws.Range("L134").Value = 17834.1432
$ws.Range("M134").Value = -2687.2761
$ws.Range("N134").Value = -22904.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 547.3333
$ws.Range("I5").Value = 383.0909
$ws.Range("J5").Value = 660.25
$ws.Range("K5").Value = 1149.2727
$ws.Range("L5").Value = 1980.75
$ws.Range("M5").Value = -1037.2727
$ws.Range("N5").Value = -2204.75
$ws.Range("H68").Value = 916.26666
$ws.Range("I68").Value = 761.93335
$ws.Range("J68").Value = 993.43335
$ws.Range("K68").Value = 2285.80005
$ws.Range("L68").Value = 2980.30005
$ws.Range("M68").Value = -1474.80005
$ws.Range("N68").Value = -4602.30005
$ws.Range("H71").Value = 916.26666
$ws.Range("I71").Value = 761.93335
$ws.Range("J71").Value = 993.43335
$ws.Range("K71").Value = 6857.40015
$ws.Range("L71").Value = 8940.900149999999
$ws.Range("M71").Value = -2801.40015
$ws.Range("N71").Value = -17052.90015
$ws.Range("H107").Value = 43478900
$ws.Range("I107").Value = 205.84616
$ws.Range("J107").Value = 100001200
$ws.Range("K107").Value = 617.5384799999999
$ws.Range("L107").Value = 300003600
$ws.Range("M107").Value = 1302.46152
$ws.Range("N107").Value = -300007440
$ws.Range("H129").Value = 7700.3125
$ws.Range("J129").Value = 25686.5
$ws.Range("L129").Value = 77059.5
$ws.Range("N129").Value = -87059.5
$ws.Range("H135").Value = 547.3333
$ws.Range("I135").Value = 383.0909
$ws.Range("J135").Value = 660.25
$ws.Range("K135").Value = 3447.8181
$ws.Range("L135").Value = 5942.25
$ws.Range("M135").Value = -912.8181
$ws.Range("N135").Value = -11012.25
$ws.Range("H141").Value = 5560264.5
$ws.Range("I141").Value = 8335396.5
$ws.Range("K141").Value = 25006189.5
$ws.Range("M141").Value = -25001009.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 9986.5
$ws.Range("J34").Value = 9986.5
$ws.Range("L34").Value = 9986.5
$ws.Range("N34").Value = -10522.5
$ws.Range("H76").Value = 9986.5
$ws.Range("J76").Value = 9986.5
$ws.Range("L76").Value = 9986.5
$ws.Range("N76").Value = -10616.5
$ws.Range("H79").Value = 9986.5
$ws.Range("J79").Value = 9986.5
$ws.Range("L79").Value = 9986.5
$ws.Range("N79").Value = -12170.5
$ws.Range("H102").Value = 2018.5
$ws.Range("I102").Value = 2150.9092
$ws.Range("K102").Value = 2150.9092
$ws.Range("M102").Value = -528.9092000000001
$ws.Range("H122").Value = 3750.7
$ws.Range("I122").Value = 4215.2856
$ws.Range("J122").Value = 2666.6667
$ws.Range("K122").Value = 12645.8568
$ws.Range("L122").Value = 8000.000100000001
$ws.Range("M122").Value = -10195.8568
$ws.Range("N122").Value = -12900.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 5613.409
$ws.Range("I93").Value = 6839.1763
$ws.Range("J93").Value = 1445.8
$ws.Range("K93").Value = 6839.1763
$ws.Range("L93").Value = 1445.8
$ws.Range("M93").Value = -5591.1763
$ws.Range("N93").Value = -3941.8
$ws.Range("H122").Value = 2983.0417
$ws.Range("I122").Value = 2535.2307
$ws.Range("J122").Value = 3512.2727
$ws.Range("K122").Value = 7605.6921
$ws.Range("L122").Value = 10536.8181
$ws.Range("M122").Value = -5155.6921
$ws.Range("N122").Value = -15436.8181

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3177.7334
$ws.Range("I122").Value = 2478
$ws.Range("K122").Value = 7434
$ws.Range("M122").Value = -4984
